$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$lines = @(
  "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)",
  "OPTIONAL MATCH (samp:sample)-->(c)",
  "OPTIONAL MATCH (diag:diagnosis)-->(c)",
  "OPTIONAL MATCH (f:file)-[*]->(c)",
  "OPTIONAL MATCH (sf:file)-->(s)",
  "WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p",
  "WHERE demo.breed IN ['Doberman Pinscher']",
  "RETURN  ",
  "    count(distinct p) AS Programs,",
  "    count(distinct s) AS Studies,",
  "    count(distinct c) AS Cases,",
  "    count(distinct samp) AS Samples,",
  "    count(distinct f) AS " + [char]96 + "Case Files" + [char]96 + ",",
  "    count(distinct sf) AS " + [char]96 + "Study Files" + [char]96
)
$newQuery = [string]::Join([char]10, $lines)

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Adjust row heights (auto-shrink after shorter query text)
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# Update sheet view: scroll back to top (removes the old A5 top-left anchor)
# and move the selection to B4:B5
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4:B5").Select()

